$d = $word.ActiveDocument

# 1. "francs" + " CFA)" merge -> "francs CFA)" (remove space-run split)
$d.Content.Find.Execute("francs CFA)", $true, $false, $false, $false, $false, $true, 1, $false, "francs CFA)", 2)

# 2. "N° " -> "N°" (trailing space removed before the bookmark)
$d.Content.Find.Execute("N° ", $true, $false, $false, $false, $false, $true, 1, $false, "N°", 2)

# 3. Date change
$d.Content.Find.Execute("Fait à Libreville, le 09 décembre 2024.", $true, $false, $false, $false, $false, $true, 1, $false, "Fait à Libreville, le 12 décembre 2024.", 2)
